# Adds the new "Sensitivity Labels" settings rows to the
# "Request Settings" sheet / Site_Request_Settings table, matching the
# SharePoint List items.xlsx update from the "sensitivitylabels support"
# commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Request Settings")

# --- Append the three new setting rows -----------------------------------
# Row 11: EnableSensitivityLabels
$ws.Range("A11").Value = "EnableSensitivityLabels"
$ws.Range("B11").Value = $false
$ws.Range("C11").Value = "Enable the Sensitivity Label functionality. "

# Row 12: DefaultSensitivityLabel
$ws.Range("A12").Value = "DefaultSensitivityLabel"
$ws.Range("C12").Value = "The default Sensitivity Label to display to users in the app. This must be a valid label id from the IP labels list."

# Row 13: RequireSensitivityLabel
$ws.Range("A13").Value = "RequireSensitivityLabel"
$ws.Range("B13").Value = $false
$ws.Range("C13").Value = "Require the user to select a Sensitivity Label in the app."

# --- Grow the query table / ListObject to cover the new rows -------------
$lo = $ws.ListObjects.Item(1)
$null = $lo.Resize($ws.Range("A1:C13"))

# --- Keep the sheet's hidden ExternalData_1 name in sync with the table --
$extName = $ws.Names.Item(1)
$extName.RefersTo = "='Request Settings'!`$A`$1:`$C`$13"

# --- Match the new selection left behind in the saved workbook -----------
$null = $ws.Range("B19").Select()
